# Scheduled-runner update: refresh Leve profit calculations (currentAveragePrice,
# NQ/HQ prices and derived profit columns) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets with freshly pulled market-board data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1551.6842
$ws.Range("I15").Value = 1551.6842
$ws.Range("K15").Value = 4655.0526
$ws.Range("M15").Value = -4486.0526

$ws.Range("H29").Value = 874
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H58").Value = 2647.6667
$ws.Range("J58").Value = 5000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15300

$ws.Range("H96").Value = 2011
$ws.Range("I96").Value = 2043.8334
$ws.Range("K96").Value = 6131.5002
$ws.Range("M96").Value = -4758.5002

$ws.Range("H100").Value = 8320672.5
$ws.Range("I100").Value = 22262
$ws.Range("K100").Value = 22262
$ws.Range("M100").Value = -21721

$ws.Range("H106").Value = 5283.8945
$ws.Range("I106").Value = 4678.2144
$ws.Range("K106").Value = 4678.2144
$ws.Range("M106").Value = -4047.2144

$ws.Range("H132").Value = 5669.2104
$ws.Range("I132").Value = 5669.2104
$ws.Range("K132").Value = 17007.6312
$ws.Range("M132").Value = -14477.6312

$ws.Range("H138").Value = 5359.1577
$ws.Range("I138").Value = 1233.8667
$ws.Range("J138").Value = 6832.476
$ws.Range("K138").Value = 3701.6001
$ws.Range("L138").Value = 20497.428
$ws.Range("M138").Value = 1438.3999
$ws.Range("N138").Value = -30777.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4334.75
$ws.Range("I61").Value = 4334.75
$ws.Range("K61").Value = 4334.75
$ws.Range("M61").Value = -4122.75

$ws.Range("H97").Value = 14225.777
$ws.Range("I97").Value = 19207.834
$ws.Range("J97").Value = 4261.6665
$ws.Range("K97").Value = 19207.834
$ws.Range("L97").Value = 4261.6665
$ws.Range("M97").Value = -18711.834
$ws.Range("N97").Value = -5253.6665

$ws.Range("H102").Value = 9172.223
$ws.Range("I102").Value = 9750
$ws.Range("J102").Value = 9100
$ws.Range("K102").Value = 9750
$ws.Range("L102").Value = 9100
$ws.Range("M102").Value = -8128
$ws.Range("N102").Value = -12344

$ws.Range("H122").Value = 1172720.6
$ws.Range("J122").Value = 1405764.9
$ws.Range("L122").Value = 4217294.699999999
$ws.Range("N122").Value = -4222194.699999999

$ws.Range("H132").Value = 3958.9443
$ws.Range("I132").Value = 3156.1667
$ws.Range("J132").Value = 5564.5
$ws.Range("K132").Value = 9468.500100000001
$ws.Range("L132").Value = 16693.5
$ws.Range("M132").Value = -6938.500100000001
$ws.Range("N132").Value = -21753.5

$ws.Range("H136").Value = 4334.75
$ws.Range("I136").Value = 4334.75
$ws.Range("K136").Value = 13004.25
$ws.Range("M136").Value = -10454.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 7250
$ws.Range("J10").Value = 7250
$ws.Range("L10").Value = 7250
$ws.Range("N10").Value = -7530

$ws.Range("H94").Value = 2319
$ws.Range("I94").Value = 2577.3333
$ws.Range("K94").Value = 2577.3333
$ws.Range("M94").Value = -2126.3333

$ws.Range("H99").Value = 45789.375
$ws.Range("I99").Value = 51116.43
$ws.Range("K99").Value = 51116.43
$ws.Range("M99").Value = -49618.43

$ws.Range("H105").Value = 74078.36
$ws.Range("I105").Value = 101795.3
$ws.Range("J105").Value = 4786
$ws.Range("K105").Value = 101795.3
$ws.Range("L105").Value = 4786
$ws.Range("M105").Value = -100048.3
$ws.Range("N105").Value = -8280

$ws.Range("H134").Value = 4898.65
$ws.Range("I134").Value = 5176.567
$ws.Range("K134").Value = 15529.701
$ws.Range("M134").Value = -12994.701

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 20009
$ws.Range("I23").Value = 20009
$ws.Range("K23").Value = 20009
$ws.Range("M23").Value = -19769

$ws.Range("H27").Value = 20009
$ws.Range("I27").Value = 20009
$ws.Range("K27").Value = 20009
$ws.Range("M27").Value = -19817

$ws.Range("H105").Value = 193417.64
$ws.Range("I105").Value = 301085.16
$ws.Range("K105").Value = 301085.16
$ws.Range("M105").Value = -299338.16

$ws.Range("H122").Value = 2130.75
$ws.Range("I122").Value = 2056.9
$ws.Range("K122").Value = 6170.700000000001
$ws.Range("M122").Value = -3720.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2996
$ws.Range("I3").Value = 2996
$ws.Range("K3").Value = 8988
$ws.Range("M3").Value = -8876

$ws.Range("H75").Value = 1338.4286
$ws.Range("I75").Value = 490
$ws.Range("J75").Value = 1974.75
$ws.Range("K75").Value = 1470
$ws.Range("L75").Value = 5924.25
$ws.Range("M75").Value = -472
$ws.Range("N75").Value = -7920.25

$ws.Range("H78").Value = 1338.4286
$ws.Range("I78").Value = 490
$ws.Range("J78").Value = 1974.75
$ws.Range("K78").Value = 4410
$ws.Range("L78").Value = 17772.75
$ws.Range("M78").Value = 582
$ws.Range("N78").Value = -27756.75

$ws.Range("H118").Value = 2731.6667
$ws.Range("I118").Value = 764.3333
$ws.Range("K118").Value = 2292.9999
$ws.Range("M118").Value = -1049.9999

$ws.Range("H132").Value = 47593
$ws.Range("I132").Value = 1470.25
$ws.Range("J132").Value = 73948.86
$ws.Range("K132").Value = 13232.25
$ws.Range("L132").Value = 665539.74
$ws.Range("M132").Value = -10702.25
$ws.Range("N132").Value = -670599.74

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 9434.559999999999
$ws.Range("I97").Value = 11139.2
$ws.Range("J97").Value = 2616
$ws.Range("K97").Value = 11139.2
$ws.Range("L97").Value = 2616
$ws.Range("M97").Value = -10643.2
$ws.Range("N97").Value = -3608

$ws.Range("H122").Value = 18593.285
$ws.Range("I122").Value = 21982.455
$ws.Range("J122").Value = 6166.3335
$ws.Range("K122").Value = 65947.36500000001
$ws.Range("L122").Value = 18499.0005
$ws.Range("M122").Value = -63497.36500000001
$ws.Range("N122").Value = -23399.0005

$ws.Range("H123").Value = 18909
$ws.Range("J123").Value = 18909
$ws.Range("L123").Value = 18909
$ws.Range("N123").Value = -23809

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1098.5
$ws.Range("I16").Value = 1098.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1098.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -928.5
$ws.Range("N16").ClearContents()

$ws.Range("H46").Value = 2168.8696
$ws.Range("I46").Value = 1390.9166
$ws.Range("J46").Value = 3017.5454
$ws.Range("K46").Value = 1390.9166
$ws.Range("M46").Value = -1202.9166
$ws.Range("N46").Value = -3393.5454

$ws.Range("H93").Value = 4279.3335
$ws.Range("I93").Value = 4314.5
$ws.Range("K93").Value = 4314.5
$ws.Range("M93").Value = -3066.5

$ws.Range("H100").Value = 5809.467
$ws.Range("I100").Value = 6261.8335
$ws.Range("K100").Value = 6261.8335
$ws.Range("M100").Value = -5720.8335

$ws.Range("H122").Value = 4156.316
$ws.Range("I122").Value = 3470.0908
$ws.Range("K122").Value = 10410.2724
$ws.Range("M122").Value = -7960.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 28126.455
$ws.Range("I81").Value = 35027.285
$ws.Range("J81").Value = 16050
$ws.Range("K81").Value = 70054.57000000001
$ws.Range("L81").Value = 32100
$ws.Range("M81").Value = -68993.57000000001
$ws.Range("N81").Value = -34222

$ws.Range("H84").Value = 28126.455
$ws.Range("I84").Value = 35027.285
$ws.Range("J84").Value = 16050
$ws.Range("K84").Value = 350272.85
$ws.Range("L84").Value = 160500
$ws.Range("M84").Value = -344968.85
$ws.Range("N84").Value = -171108

$ws.Range("H96").Value = 4990
$ws.Range("I96").Value = 4990
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 4990
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -3617
$ws.Range("N96").ClearContents()

$ws.Range("H100").Value = 66668.5
$ws.Range("I100").Value = 37891.5
$ws.Range("J100").Value = 152999.5
$ws.Range("K100").Value = 75783
$ws.Range("L100").Value = 305999
$ws.Range("M100").Value = -75242
$ws.Range("N100").Value = -307081

$ws.Range("H136").Value = 3012.7407
$ws.Range("I136").Value = 2232.6086
$ws.Range("J136").Value = 7498.5
$ws.Range("K136").Value = 6697.825800000001
$ws.Range("L136").Value = 22495.5
$ws.Range("M136").Value = -4147.825800000001
$ws.Range("N136").Value = -27595.5
